$p = $ppt.ActivePresentation

# --- 1. Table on slide 5: switch to the new built-in table style ---
$tableSlide = $p.Slides.Item(5)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{820AE70F-8B0D-4158-92C4-98C397B7242F}")

# --- 2. Re-point the presentation's theme colour scheme from the old
#        "Integral / Red Violet" palette to the standard "Office Theme"
#        palette (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink). ---
$theme = $p.SlideMaster.Theme
$clrScheme = $theme.ThemeColorScheme

$clrScheme.Item(1).RGB  = 0x000000   # dk1
$clrScheme.Item(2).RGB  = 0xFFFFFF   # lt1
$clrScheme.Item(3).RGB  = 0x6A5444   # dk2      (BGR order: 44546A)
$clrScheme.Item(4).RGB  = 0xE6E6E7   # lt2      (BGR order: E7E6E6)
$clrScheme.Item(5).RGB  = 0xD59B5B   # accent1  (BGR order: 5B9BD5)
$clrScheme.Item(6).RGB  = 0x317DED   # accent2  (BGR order: ED7D31)
$clrScheme.Item(7).RGB  = 0xA5A5A5   # accent3  (BGR order: A5A5A5)
$clrScheme.Item(8).RGB  = 0x00C0FF   # accent4  (BGR order: FFC000)
$clrScheme.Item(9).RGB  = 0xC47244   # accent5  (BGR order: 4472C4)
$clrScheme.Item(10).RGB = 0x47AD70   # accent6  (BGR order: 70AD47)
$clrScheme.Item(11).RGB = 0xC16305   # hlink    (BGR order: 0563C1)
$clrScheme.Item(12).RGB = 0x724F95   # folHlink (BGR order: 954F72)
